# Retro Paradise deck edit:
#  1. Fix "Tech Involved" slide text: "SQLite or " -> "SQLite/"
#  2. Insert a brand-new "User Stories" slide at position 4 (between
#     "Tech Involved" and "Context Diagram"), using the same
#     "Title and Content" layout as the surrounding slides.

$p = $ppt.ActivePresentation

# --- 1. Tech Involved slide: tidy "SQLite or " -> "SQLite/" -------------
$techSlide = $p.Slides.Item(3)
$techBody  = $techSlide.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $techBody.Paragraphs().Count; $i++) {
    $para = $techBody.Paragraphs($i, 1)
    if ($para.Text -like "SQLite or *") {
        $run = $para.Characters(1, 10)
        $run.Text = "SQLite/"
    }
}

# --- 2. Insert new "User Stories" slide at index 4 ----------------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"
$newSlide = $p.Slides.AddSlide(4, $layout)

# Title
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "User Stories"

# Body content placeholder
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "As a customer`r" + `
             "   I want an online catalog`r" + `
             "   So I can browse store products `r" + `
             "`r" + `
             " As an Admin`r" + `
             "    I want a database`r" + `
             "    So I can add products and check stock"

# Apply the drop-shadow text effect to the "customer story" paragraphs
# (paragraphs 1-3) exactly like the source deck.
for ($i = 1; $i -le 3; $i++) {
    $para = $body.Paragraphs($i, 1)
    $para.Font.Shadow = $true
}

# Remove the bullet / auto-indent on the wrapped continuation lines so they
# read as plain indented text under "As a customer" / "As an Admin".
$para2 = $body.Paragraphs(2, 1)
$para2.ParagraphFormat.Bullet.Visible = $false

$para3 = $body.Paragraphs(3, 1)
$para3.ParagraphFormat.Bullet.Visible = $false

$para4 = $body.Paragraphs(4, 1)
$para4.ParagraphFormat.Bullet.Visible = $false

$para6 = $body.Paragraphs(6, 1)
$para6.ParagraphFormat.Bullet.Visible = $false

$para7 = $body.Paragraphs(7, 1)
$para7.ParagraphFormat.Bullet.Visible = $false

# Split paragraphs 2, 3 and 7 into the same two runs as the source deck
# (the lead-in word gets its own run before the rest of the sentence).
$para2.Characters(1, 5).Text = $para2.Characters(1, 5).Text
$para3.Characters(1, 6).Text = $para3.Characters(1, 6).Text
$para7.Characters(1, 30).Text = $para7.Characters(1, 30).Text
